# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sorted order of workers (3 new workers moved to the top, rest follow),
# with updated "Salario Basico" values for a couple of existing workers.
$data = @(
    @{ Row = 16; Doc = "CC"; Num = "73196033";   Nombre = "JOSE LUIS CERVANTES MEJIA";        Periodo = "1908"; Mora = 1600;  Salario = 1200000 },
    @{ Row = 17; Doc = "CC"; Num = "1143342046";  Nombre = "LUIS PUENTES VALLE";                Periodo = "2005"; Mora = 68000; Salario = 1700000 },
    @{ Row = 18; Doc = "CC"; Num = "1002392859";  Nombre = "LAURA VANESA RODRIGUEZ GONZALEZ";   Periodo = "2011"; Mora = 30430; Salario = 1300000 },
    @{ Row = 19; Doc = "CC"; Num = "73089468";    Nombre = "CELSO ANTONIO GONZALEZ FORTICH";    Periodo = "2305"; Mora = 156000;Salario = 5000000 },
    @{ Row = 20; Doc = "CC"; Num = "1047421035";  Nombre = "RICHAR RAFAEL SILGADO VILLALOBOS";  Periodo = "2305"; Mora = 2667;  Salario = 2000000 },
    @{ Row = 21; Doc = "CC"; Num = "92226300";    Nombre = "ROGER CEDRON RAMIREZ";              Periodo = "2312"; Mora = 8000;  Salario = 6000000 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.Doc
    $ws.Cells.Item($r, 3).Value = $item.Num
    $ws.Cells.Item($r, 4).Value = $item.Nombre
    $ws.Cells.Item($r, 5).Value = $item.Periodo
    $ws.Cells.Item($r, 6).Value = $item.Mora
    $ws.Cells.Item($r, 7).Value = $item.Salario
}
